$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update the first time-sheet entry (date + time in/out)
$ws.Range("A2").Value = 45563
$ws.Range("B2").Value = 0.625
$ws.Range("C2").Value = 0.70833333333333337

# Row 3: new time-sheet entry — date auto-filled via TODAY(), plus time in/out
$ws.Range("A3").Formula = "=TODAY()"
$ws.Range("B3").Value = 0.54166666666666663
$ws.Range("C3").Value = 0.60416666666666663

# Move the active selection to A4
$ws.Range("A4").Select() | Out-Null
